$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.239.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +6.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.282.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.72"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "632.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.415"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +15.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.712"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.998"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.278.96"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.589"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000266"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.85%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.21"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.887.55"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.768.91"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.37"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.306.78"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.31"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.83%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000192"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +50.74%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Polkadot"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.29"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.38"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.21"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.514.49"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.56"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.182"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.82"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "562.65"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.14"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.68"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +26.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.93"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.29"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.85%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.44"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.397"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.01"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.53"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "181.25"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.79"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.29"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.634"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.20"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.38%  "
